$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conditions")

# Remove the "saveat[]" column entirely (column E), shifting tspan/observables[] left.
$ws.Range("E1").EntireColumn.Delete()

# tspan (now column E) for the "withdata2" row (row 3) gets a value.
$ws.Range("E3").Value = 200

# Match the author's final selection (cell E3 - the one just edited).
$ws.Range("E3").Select()
